# Add files via upload
# - Adds a "조점수" (group score) label at G2 and moves the per-row
#   "K{row}+AVERAGEIF(...)" formula up from G14:G23 to G3:G12 (the label
#   that used to live at G13 moves to G2, and G13/G14:G23 become empty).
# - Adds a new summary row 24 with a "평균" (average) label in B24 and
#   an AVERAGE(C3:C23) formula in C24.
# - Moves the active selection to F14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Group-score header label, now shown at G2 (previously at G13). Copy the
# number format from the old label cell so G2 keeps the same style index
# instead of minting a new one, then set the text.
$ws.Range("G13").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = "조점수"

# Move the per-student group-score formulas from rows 14-23 up to rows 3-12.
for ($i = 0; $i -le 9; $i++) {
    $srcRow = 14 + $i
    $dstRow = 3 + $i
    $ws.Range("G$dstRow").Formula = "=K$srcRow+AVERAGEIF(L${srcRow}:AF$srcRow,`">3.0`")"
    $ws.Range("G$srcRow").Clear()
}

# The old label cell at G13 is no longer needed.
$ws.Range("G13").Clear()

# New average row.
$ws.Range("B24").Value = "평균"
$ws.Range("C24").Formula = "=AVERAGE(C3:C23)"

# Match the saved selection state.
$ws.Range("F14").Select()
